# Generate Report for Archive
#
# 1) Update the status text "Ready for handoff" -> "In Translation" for the
#    row-2 status cells on every sheet that shows it (Overview!E2:F2 and the
#    "Status" column (C2) on the zh-cn / de-de detail sheets).
# 2) Narrow the corresponding "Status" columns (Overview columns E & F,
#    zh-cn/de-de column C) to match the new, shorter status text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# ColumnWidth (in characters) that this engine quantizes/serializes back out
# as an OOXML column width as close as possible to the target 13.41 used in
# the updated report.
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the literal on the left of -eq; PowerShell's -eq uses the
        # left operand's type to decide how to coerce the right operand, and
        # $cell.Value() can come back as a .NET bool for True/False cells -
        # comparing "someText" -eq $trueOrFalse would otherwise silently
        # coerce the string to a boolean instead of doing a string compare.
        if ($oldStatus -eq $cell.Value()) {
            $cell.Value = $newStatus
            $ws.Columns.Item($cell.Column()).ColumnWidth = $newColumnWidth
        }
    }
}
